$d = $word.ActiveDocument

# Helper: replace text while preserving a preceding empty run (<w:r/>) that
# would otherwise get coalesced into the freshly written run because both
# share identical (default) formatting. Toggling Bold on/off around the
# text assignment forces the engine to keep the runs distinct.
function Replace-KeepEmptyRun($doc, [string]$old, [string]$new) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = 1
        $rng.Text = $new
        $rng.Font.Bold = 0
    }
}

# 1: Main H1 title (no preceding empty run - plain replace is safe)
$d.Content.Find.Execute("Play Maximus Payus for Free: Review and Details", $true, $false, $false, $false, $false, $true, 1, $false, "Play Maximus Payus for Free", 2)

# 2: Bullet under "What we like"
Replace-KeepEmptyRun $d "Free spins with increasing multipliers" "Thematic and visually appealing graphics"

# 3: Bullet under "What we like"
Replace-KeepEmptyRun $d "Random bonuses in regular play" "Free spins and random bonuses"

# 4: Bullet under "What we like"
Replace-KeepEmptyRun $d "Affordable betting options" "Accessible betting options for all players"

# 5: Bullet under "What we don't like"
Replace-KeepEmptyRun $d "High volatility may not be suitable for all players" "High volatility may result in less frequent wins"

# 6: Bullet under "What we don't like"
Replace-KeepEmptyRun $d "Graphics may not appeal to those looking for a more realistic theme" "No progressive jackpot feature"

# 7 (same replacement as step 1, which already updated both occurrences of
# the title/heading text via wdReplaceAll, including the bold run near the
# end) - nothing further required here.

# 8: Meta description italic paragraph (differing formatting from its
# preceding empty run - plain replace is safe)
$d.Content.Find.Execute("Read our review of Maximus Payus, a high-volatility slot game with 1,024 ways to win and free spins with increasing multipliers. Play for free now.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Maximus Payus slot game and play for free. Enjoy 1,024 ways to win and exciting bonuses.", 2)
